$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New cell G11/H11: copy the format from G10/H10 (same little "rua" legend block)
# and fill in the new "rua2" entry.
$ws.Range("G10:H10").Copy()
$ws.Range("G11").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("G11").Value = 2
$ws.Range("H11").Value = "rua2"

# --- L11 changed from 3 to 1
$ws.Range("L11").Value = 1

# --- G18 changed from 4 to 2
$ws.Range("G18").Value = 2

# --- New underlined, blank cell J7
$ws.Range("J7").Font.Underline = $true

# --- Sheet view: selection now on L11 (the workbook was also scrolled so that
# column E is the leftmost visible column, but that pure view/scroll position
# isn't part of the bridged object model here, so only the selection is settable)
$ws.Range("L11").Select()

# --- Page setup: Letter/A4-ish paper size 9 (A4), portrait orientation
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$wb.Save()
